$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: add new mail-log row 4 ---
$logs.Range("A4").Value = "Offerte voor zakelijke samenwerking"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("C4").Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$logs.Range("D4").Value = "Bestelling"
$logs.Range("F4").Value = "2025-06-18 09:00:12"
$logs.Range("G4").Value = "Nee"

# --- Extend conditional formatting ranges to include row 4 ---
$logs.Range("D2:D3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D4"))
$logs.Range("G2:G3").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G4"))

# --- Dashboard sheet: add new category total row 4 ---
$dash.Range("A4").Value = "Bestelling"
$dash.Range("B4").Value = 1

# --- Extend chart series ranges to include the new Dashboard row ---
$co = $dash.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$4,Dashboard!`$B`$2:`$B`$4,1)"
